$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '98.804.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.56%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.313.27'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.62%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '255.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '626.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.50%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.46'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +21.78%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.416'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.05%  '

$ws.Range('E9').Value = '  +0.01%  '

$ws.Range('E10').Value = '  +24.64%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.312.54'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.54%  '

$ws.Range('E12').Value = '  +2.83%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '41.38'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +15.30%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '98.479.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000253'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.60%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.941.91'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.24%  '

$ws.Range('E17').Value = '  -1.63%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.315.93'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.65%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.30%  '

$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.08%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.95%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '486.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.86%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.36%  '

$ws.Range('E24').Value = '  -2.48%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.81'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.24%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.344'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +35.59%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '89.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.40%  '

$ws.Range('E28').Value = '  +1.13%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.494.04'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.44%  '

$ws.Range('E30').Value = '  +20.53%  '

$ws.Range('E31').Value = '  -0.10%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.191'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.59%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +15.34%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.07%  '

$ws.Range('E35').Value = '  +2.47%  '

$ws.Range('E36').Value = '  +7.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.151'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.94%  '

$ws.Range('E38').Value = '  +0.04%  '

$ws.Range('E39').Value = '  +0.80%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '498.28'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.49%  '

$ws.Range('B41').Value = 'MantraDAO'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.94'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.61%  '

$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.74'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.24%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.58%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.791'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.21%  '

$ws.Range('E45').Value = '  +0.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.18'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '160.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.54%  '

$ws.Range('E48').Value = '  +1.39%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.859'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.99%  '

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.76%  '

$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.95%  '
